$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '58.169.79'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.59%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.140.74'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  -0.01%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '534.18'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '138.85'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("E7").Value = '  +0.02%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '3.139.15'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("E9").Value = '  +4.50%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '7.30'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.59%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.107'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.49%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.414'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +4.59%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '3.680.53'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.13%  '

$ws.Range("E14").Value = '  +1.56%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '25.67'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.84%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.0000164'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '58.239.95'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.47%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.144.79'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.44%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '6.06'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.88%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '12.71'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.49%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '8.17'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.75%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '359.51'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +1.57%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '69.09'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.506'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("E28").Value = '  -3.55%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '7.32'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.29%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '6.17'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("E31").Value = '  -0.20%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '21.51'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.55%  '

$ws.Range("E33").Value = '  +1.15%  '

$ws.Range("E34").Value = '  -2.73%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '159.37'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.32%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.08'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -1.39%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '26.05'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("E38").Value = '  +0.92%  '

$ws.Range("E39").Value = '  +4.85%  '

$ws.Range("E40").Value = '  +0.09%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.508.51'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +8.50%  '

$ws.Range("E42").Value = '  -0.02%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '4.01'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -4.26%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '37.42'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +2.41%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.182.70'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("E49").Value = '  +0.51%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '19.85'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -2.79%  '

$ws.Range("E51").Value = '  -4.11%  '

